$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.452.22'
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').Value = '1.583.60'
$ws.Range('E3').Value = '  -0.24%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '213.62'
$ws.Range('E5').Value = '  +0.36%  '

$ws.Range('E6').Value = '  +0.29%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = '44.34'
$ws.Range('E8').Value = '  +0.19%  '

$ws.Range('D9').Value = '24.04'
$ws.Range('E9').Value = '  -0.89%  '

$ws.Range('E10').Value = '  -1.81%  '

$ws.Range('E11').Value = '  -1.48%  '

$ws.Range('E12').Value = '  +0.88%  '

$ws.Range('E13').Value = '  -0.40%  '

$ws.Range('D14').Value = '1.593.31'
$ws.Range('E14').Value = '  +0.55%  '

$ws.Range('E15').Value = '  -1.00%  '

$ws.Range('D16').Value = '0.521'
$ws.Range('E16').Value = '  -1.44%  '

$ws.Range('D17').Value = '28.463.53'
$ws.Range('E17').Value = '  -0.16%  '

$ws.Range('D18').Value = '62.23'
$ws.Range('E18').Value = '  -1.38%  '

$ws.Range('D19').Value = '230.84'
$ws.Range('E19').Value = '  -0.26%  '

$ws.Range('D20').Value = '7.44'
$ws.Range('E20').Value = '  -0.39%  '

$ws.Range('D21').Value = '0.0₃0691'
$ws.Range('E21').Value = '  -2.18%  '

$ws.Range('E22').Value = '  +0.17%  '

$ws.Range('E23').Value = '  -3.12%  '

$ws.Range('E24').Value = '  -1.90%  '

$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  +4.12%  '

$ws.Range('D26').Value = '152.19'
$ws.Range('E26').Value = '  +0.30%  '

$ws.Range('D27').Value = '15.03'
$ws.Range('E27').Value = '  -1.37%  '

$ws.Range('E28').Value = '  -1.58%  '

$ws.Range('E29').Value = '  -2.06%  '

$ws.Range('E31').Value = '  +2.44%  '

$ws.Range('E32').Value = '  -1.59%  '

$ws.Range('E33').Value = '  -1.23%  '

$ws.Range('E34').Value = '  -2.03%  '

$ws.Range('D35').Value = '1.399.36'
$ws.Range('E35').Value = '  -0.36%  '

$ws.Range('E36').Value = '  +6.08%  '

$ws.Range('E37').Value = '  -4.15%  '

$ws.Range('E38').Value = '  +0.57%  '

$ws.Range('E39').Value = '  +1.73%  '

$ws.Range('D40').Value = '0.0165'
$ws.Range('E40').Value = '  -0.80%  '

$ws.Range('D41').Value = '0.522'
$ws.Range('E41').Value = '  -3.59%  '

$ws.Range('E42').Value = '  +0.06%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '0.791'
$ws.Range('E43').Value = '  -2.69%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '1.89'
$ws.Range('E44').Value = '  +1.27%  '

$ws.Range('D45').Value = '0.0464'
$ws.Range('E45').Value = '  -0.56%  '

$ws.Range('D46').Value = '5.45'
$ws.Range('E46').Value = '  -3.23%  '

$ws.Range('E47').Value = '  -2.31%  '

$ws.Range('D48').Value = '63.26'
$ws.Range('E48').Value = '  +0.41%  '

$ws.Range('D49').Value = '1.719.42'
$ws.Range('E49').Value = '  -0.37%  '

$ws.Range('D50').Value = '86.64'
$ws.Range('E50').Value = '  -0.72%  '

$ws.Range('E51').Value = '  -0.92%  '
